# This script re-applies a scheduled market-data refresh to the per-job
# "Leve" profit-tracking sheets (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR).
# For each affected row, the live market price columns
# (H currentAveragePrice, I currentAveragePriceNQ, J currentAveragePriceHQ,
#  K LevePriceNQ, L LevePriceHQ, M LeveProfitNQ, N LeveProfitHQ) are updated
# to the freshly-fetched values.
$wb = $excel.ActiveWorkbook

# Sheet "ALC", row 100 (Asking for a Friend | Beetle Glue)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 22224860
$ws.Range("I100").Value = 47620610
$ws.Range("J100").Value = 3583.125
$ws.Range("K100").Value = 47620610
$ws.Range("L100").Value = 3583.125
$ws.Range("M100").Value = -47620069
$ws.Range("N100").Value = -4665.125

# Sheet "ARM", row 32 (Ingot We Trust | Steel Ingot)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4451.12
$ws.Range("I32").Value = 4451.12
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 4451.12
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -4164.12
$ws.Range("N32").ClearContents()

# Sheet "BSM", row 107 (The Gold Experience | Deepgold Nugget)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3845.2058
$ws.Range("I107").Value = 4333.393
$ws.Range("K107").Value = 4333.393
$ws.Range("M107").Value = -2413.393

# Sheet "BSM", row 134 (Ruthenium Supremium | Ruthenium Ingot)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 926
$ws.Range("I134").Value = 903.1667
$ws.Range("K134").Value = 2709.5001
$ws.Range("M134").Value = -174.5001000000002

# Sheet "CRP", row 86 (Birch, Please | Birch Lumber)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 3675.889
$ws.Range("I86").Value = 3136.1667
$ws.Range("K86").Value = 3136.1667
$ws.Range("M86").Value = -2013.1667

# Sheet "CRP", row 89 (Built This City on Blocks and Soul (L) | Birch Lumber)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 3675.889
$ws.Range("I89").Value = 3136.1667
$ws.Range("K89").Value = 15680.8335
$ws.Range("M89").Value = -10064.8335

# Sheet "CRP", row 132 (Hull Lotta Damage | Ginseng Lumber)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1060.5
$ws.Range("I132").Value = 903.375
$ws.Range("J132").Value = 1374.75
$ws.Range("K132").Value = 2710.125
$ws.Range("L132").Value = 4124.25
$ws.Range("M132").Value = -180.125
$ws.Range("N132").Value = -9184.25

# Sheet "CUL", row 5 (What a Sap | Maple Syrup)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1141.3125
$ws.Range("I5").Value = 688.875
$ws.Range("J5").Value = 1593.75
$ws.Range("K5").Value = 2066.625
$ws.Range("L5").Value = 4781.25
$ws.Range("M5").Value = -1954.625
$ws.Range("N5").Value = -5005.25

# Sheet "CUL", row 68 (Such a Butter Face | Fermented Butter)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1092.0834
$ws.Range("I68").Value = 920.4
$ws.Range("J68").Value = 1214.7142
$ws.Range("K68").Value = 2761.2
$ws.Range("L68").Value = 3644.1426
$ws.Range("M68").Value = -1950.2
$ws.Range("N68").Value = -5266.142599999999

# Sheet "CUL", row 71 (No Margarine of Error (L) | Fermented Butter)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 1092.0834
$ws.Range("I71").Value = 920.4
$ws.Range("J71").Value = 1214.7142
$ws.Range("K71").Value = 8283.6
$ws.Range("L71").Value = 10932.4278
$ws.Range("M71").Value = -4227.6
$ws.Range("N71").Value = -19044.4278

# Sheet "CUL", row 131 (The Mountain Steeped | Tsai tou Vounou)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 6850299
$ws.Range("I131").Value = 1021.0526
$ws.Range("J131").Value = 9260230
$ws.Range("K131").Value = 3063.1578
$ws.Range("L131").Value = 27780690
$ws.Range("M131").Value = 1976.8422
$ws.Range("N131").Value = -27790770

# Sheet "CUL", row 132 (More Mezcal | Cooking Mezcal)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1938.5
$ws.Range("I132").Value = 1251.1765
$ws.Range("J132").Value = 2494.9048
$ws.Range("K132").Value = 11260.5885
$ws.Range("L132").Value = 22454.1432
$ws.Range("M132").Value = -8730.5885
$ws.Range("N132").Value = -27514.1432

# Sheet "CUL", row 135 (Not-so-secret Ingredient | Royal Maple Syrup)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1141.3125
$ws.Range("I135").Value = 688.875
$ws.Range("J135").Value = 1593.75
$ws.Range("K135").Value = 6199.875
$ws.Range("L135").Value = 14343.75
$ws.Range("M135").Value = -3664.875
$ws.Range("N135").Value = -19413.75

# Sheet "CUL", row 137 (Creative Chocolate | Gateau au Chocolat)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 4494.6343
$ws.Range("I137").Value = 4722.7407
$ws.Range("J137").Value = 4054.7144
$ws.Range("K137").Value = 14168.2221
$ws.Range("L137").Value = 12164.1432
$ws.Range("M137").Value = -9068.222100000001
$ws.Range("N137").Value = -22364.1432

# Sheet "GSM", row 107 (Whetstones for the Workers | Hard Mudstone Whetstone)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 335.88
$ws.Range("I107").Value = 313.42856
$ws.Range("J107").Value = 453.75
$ws.Range("K107").Value = 313.42856
$ws.Range("L107").Value = 453.75
$ws.Range("M107").Value = 1606.57144
$ws.Range("N107").Value = -4293.75

# Sheet "GSM", row 132 (On Board for Lar | Lar Ingot)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2295.8948
$ws.Range("I132").Value = 1859.7587
$ws.Range("J132").Value = 3701.2222
$ws.Range("K132").Value = 5579.2761
$ws.Range("L132").Value = 11103.6666
$ws.Range("M132").Value = -3049.2761
$ws.Range("N132").Value = -16163.6666

# Sheet "GSM", row 140 (The Right Rod | Ra'Kaznar Rod)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H140").Value = 89889.5
$ws.Range("J140").Value = 89889.5
$ws.Range("L140").Value = 89889.5
$ws.Range("N140").Value = -100249.5

# Sheet "LTW", row 16 (Saddle Sore | Hard Leather)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 971.4286
$ws.Range("I16").Value = 966.6667
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 966.6667
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -796.6667
$ws.Range("N16").Value = -1340

# Sheet "LTW", row 122 (Hell on Leather | Gaja Leather)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2846.647
$ws.Range("I122").Value = 2504.4285
$ws.Range("J122").Value = 4443.6665
$ws.Range("K122").Value = 7513.2855
$ws.Range("L122").Value = 13330.9995
$ws.Range("M122").Value = -5063.2855
$ws.Range("N122").Value = -18230.9995

# Sheet "LTW", row 123 (Running up the Tabi | Gajaskin Tabi)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H123").Value = 10000
$ws.Range("J123").Value = 10000
$ws.Range("L123").Value = 10000
$ws.Range("N123").Value = -19800

# Sheet "LTW", row 136 (Respect for Br'aax | Br'aax Leather)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3561.0488
$ws.Range("I136").Value = 2414.3928
$ws.Range("J136").Value = 6030.769
$ws.Range("K136").Value = 7243.178400000001
$ws.Range("L136").Value = 18092.307
$ws.Range("M136").Value = -4693.178400000001
$ws.Range("N136").Value = -23192.307

# Sheet "WVR", row 81 (Where the Dragonflies, the Net Catches | Crawler Silk)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 38467180
$ws.Range("I81").Value = 1411.5555
$ws.Range("J81").Value = 125015160
$ws.Range("K81").Value = 2823.111
$ws.Range("L81").Value = 250030320
$ws.Range("M81").Value = -1762.111
$ws.Range("N81").Value = -250032442

# Sheet "WVR", row 84 (To Kill a Dragon on Nameday (L) | Crawler Silk)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 38467180
$ws.Range("I84").Value = 1411.5555
$ws.Range("J84").Value = 125015160
$ws.Range("K84").Value = 14115.555
$ws.Range("L84").Value = 1250151600
$ws.Range("M84").Value = -8811.555
$ws.Range("N84").Value = -1250162208

# Sheet "WVR", row 107 (Flax Wax | Bright Linen Yarn)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1241.8077
$ws.Range("I107").Value = 1213
$ws.Range("J107").Value = 1400.25
$ws.Range("K107").Value = 3639
$ws.Range("L107").Value = 4200.75
$ws.Range("M107").Value = -1719
$ws.Range("N107").Value = -8040.75

# Sheet "WVR", row 122 (Heavy Armoire | Dark Hempen Cloth)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2048.7036
$ws.Range("I122").Value = 1519.762
$ws.Range("J122").Value = 3900
$ws.Range("K122").Value = 4559.286
$ws.Range("L122").Value = 11700
$ws.Range("M122").Value = -2109.286
$ws.Range("N122").Value = -16600

# Sheet "WVR", row 126 (A Polished Purchase | Snow Linen)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1617.5
$ws.Range("I126").Value = 1420.091
$ws.Range("J126").Value = 3789
$ws.Range("K126").Value = 4260.272999999999
$ws.Range("L126").Value = 11367
$ws.Range("M126").Value = -1790.272999999999
$ws.Range("N126").Value = -16307
